# Generate Report for Handback
# Update the "last generated" timestamp values on each worksheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-08-20 05:06:09"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file
$wsZhCn.Range("H2").Value = "2016-08-20 05:06:00"
$wsZhCn.Range("K2").Value = "2016-08-20 05:06:28"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file
$wsDeDe.Range("H2").Value = "2016-08-20 05:06:09"
$wsDeDe.Range("K2").Value = "2016-08-20 05:06:35"
